# Apply the requested edits to the workbook:
# 1. Rename worksheets
# 2. Update Version and Date values on the Metadata sheet

$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Include from FBOE").Name = "Include #1"
$wb.Worksheets.Item("Exclude from FBOE").Name = "Exclude #2"

# --- Update metadata values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.1.0"
$meta.Range("B8").Value = "2024-10-31T20:37:15+01:00"
